$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: "This is an individual assignment." becomes a mail-merge
# greeting: "This is an individual assignment for «Surname», «First_name»,
# «ID_number»" using MERGEFIELD fields.
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$null = $p5.Range.Find.Execute("assignment.", $true, $false, $false, $false, `
    $false, $true, 1, $false, "assignment for ", 2)

$p = $d.Paragraphs(5)
$ins = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$null = $d.Fields.Add($ins, 59, "MERGEFIELD Surname", $false)

$p = $d.Paragraphs(5)
$ins = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$ins.InsertAfter(", ")

$p = $d.Paragraphs(5)
$ins = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$null = $d.Fields.Add($ins, 59, "MERGEFIELD First_name", $false)

$p = $d.Paragraphs(5)
$ins = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$ins.InsertAfter(", ")

$p = $d.Paragraphs(5)
$ins = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$null = $d.Fields.Add($ins, 59, "MERGEFIELD ID_number", $false)

# ---------------------------------------------------------------------
# Edit 2: "At the end of the quiz the user is shown a summary..." has
# its leading "At the end of the quiz the" replaced by a MERGEFIELD
# named "cover", leaving " user is shown a summary of how well they did
# in a graphical way using a gauge."
# ---------------------------------------------------------------------
$p11 = $d.Paragraphs(11)
$null = $p11.Range.Find.Execute("At the end of the quiz the", $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 2)

$p11 = $d.Paragraphs(11)
$insPoint = $d.Range($p11.Range.Start, $p11.Range.Start)
$null = $d.Fields.Add($insPoint, 59, "MERGEFIELD cover", $false)

Write-Output ("Para5: " + $d.Paragraphs(5).Range.Text)
Write-Output ("Para11: " + $d.Paragraphs(11).Range.Text)
